$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: AHC30716 / 110597F / CERT III HORTICULTURE ---
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "CERTIFICATE III IN HORTICULTURE"
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("H2").WrapText = $true
$ws.Range("M2").Value = "TAS"

# --- Row 3: AHC40416 / 110598E / CERT IV HORTICULTURE ---
$ws.Range("A3").Value = "AHC40416"
$ws.Range("B3").Value = "110598E"
$ws.Range("D3").Value = "CERTIFICATE IV IN HORTICULTURE"
$ws.Range("H3").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("H3").WrapText = $true
$ws.Range("M3").Value = "TAS"

# --- Row 4: AHC51422 / 110774E / DIPLOMA OF AGRIBUSINESS MANAGEMENT ---
$ws.Range("A4").Value = "AHC51422"
$ws.Range("B4").Value = "110774E"
$ws.Range("D4").Value = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("H4").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("H4").WrapText = $true
$ws.Range("M4").Value = "TAS"

# --- Row 5: package CERT III + CERT IV ---
$ws.Range("A5").Value = "AHC30716 / AHC40416"
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "110597F / 110598E"
$ws.Range("B5").WrapText = $true
$ws.Range("D5").Value = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
$ws.Range("D5").WrapText = $true
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("H5").WrapText = $true
$ws.Range("M5").Value = "TAS"

# --- Row 6: package CERT III + DIPLOMA ---
$ws.Range("A6").Value = "AHC30716 / AHC51422"
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Value = "110597F / 110774E"
$ws.Range("B6").WrapText = $true
$ws.Range("D6").Value = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("D6").WrapText = $true
$ws.Range("H6").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("H6").WrapText = $true
$ws.Range("M6").Value = "TAS"

# --- Row 7: package CERT IV + DIPLOMA ---
$ws.Range("A7").Value = "AHC40416 / AHC51422"
$ws.Range("A7").WrapText = $true
$ws.Range("B7").Value = "110598E / 110774E"
$ws.Range("B7").WrapText = $true
$ws.Range("D7").Value = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("D7").WrapText = $true
$ws.Range("H7").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("H7").WrapText = $true
$ws.Range("M7").Value = "TAS"

# --- Column C (department) filled across all rows after the main pass ---
$ws.Range("C2").Value = "HORTICULTURE"
$ws.Range("C3").Value = "HORTICULTURE"
$ws.Range("C4").Value = "MANAGEMENT"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("C6").Value = "PACKAGES"
$ws.Range("C7").Value = "PACKAGES"

# --- Column E (durationMin) numeric ---
$ws.Range("E2").Value = 52
$ws.Range("E3").Value = 52
$ws.Range("E4").Value = 52
$ws.Range("E5").Value = 104
$ws.Range("E6").Value = 104
$ws.Range("E7").Value = 104

# --- Column I (tuition) numeric, thousands format ---
$ws.Range("I2").Value = 8700
$ws.Range("I3").Value = 8700
$ws.Range("I4").Value = 9200
$ws.Range("I5").Value = 16700
$ws.Range("I6").Value = 17200
$ws.Range("I7").Value = 17200
$ws.Range("I2:I7").NumberFormat = "#,##0"

# --- Column J (tuitionDetail) filled across all rows, wrap + thousands format ---
$ws.Range("J2").Value = "8,500 tuition fee + 200 handling fee"
$ws.Range("J3").Value = "8,500 tuition fee + 200 handling fee"
$ws.Range("J4").Value = "9,000 tuition fee + 200 handling fee"
$ws.Range("J5").Value = "16,500 tuition fee + 200 handling fee"
$ws.Range("J6").Value = "17,000 tuition fee + 200 handling fee"
$ws.Range("J7").Value = "17,000 tuition fee + 200 handling fee"
$ws.Range("J2:J7").NumberFormat = "#,##0"
$ws.Range("J2:J7").WrapText = $true

# --- Column R (promotionValidity) filled across all rows last ---
$ws.Range("R2").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R3").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R4").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R5").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R6").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R7").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Row heights (45pt, matches the wrapped multi-line content) ---
$ws.Range("A2:R7").RowHeight = 45

# --- Selection matches the post-edit saved state ---
$ws.Range("S14").Select()
